$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('7917707','Chicken Leg Quarters','1','76.45','76.45'),
    @('3157690','Rice Jasmine','1','36.24','36.24'),
    @('4096954','Lentils - Dried Bean','1','24.60','24.60'),
    @('7524051','Hot Sauce - Franks Gal','1','60.40','60.40'),
    @('5818380','Can - Coconut Milk','2','92.59','185.18'),
    @('5833462','Jalapenos (Sliced)','4','26.49','105.96'),
    @('2366607','Liquid Egg','1','40.06','40.06'),
    @('7226918','Mozz (Sliced)','2','26.90','53.80'),
    @('8489148','Smoked Turkey (Unsliced)','6','76.86','461.16'),
    @('8768871','Juice - Lime','1','74.49','74.49'),
    @('7143223','Sausage - Vegan Patty','2','88.88','177.76'),
    @('3413107','Brie','2','27.49','54.98'),
    @('4171902','Manchego Cheese','1','79.74','79.74'),
    @('3617531','Roasted Root Vegetables','2','58.79','117.58'),
    @('4438911','Grain Blend 5 Way','3','77.69','233.07'),
    @('7133036','Chicken - Diced White & Dark','8','40.96','327.68'),
    @('6986178','Spanakopita - Mini','4','65.79','263.16'),
    @('6985683','Quiche - Mini (Asst)','4','70.95','283.80'),
    @('7213975','Potsticker - Pork & Veg','1','42.05','42.05'),
    @('7101552','Potsticker - Vegetable','1','57.45','57.45'),
    @('0017354','Hash Brown','2','23.89','47.78'),
    @('4617771','Egg Patty','10','46.94','469.40'),
    @('4828802','Heavy Cream','1','52.88','52.88'),
    @('1735372','Cucumber (English)','4','13.13','52.52'),
    @('1821537','Garlic - Fresh (Peeled)','1','74.67','74.67'),
    @('2004547','Herb - Basil (Fresh)','4','10.94','43.76'),
    @('1008143','Kale - Fresh','2','24.90','49.80'),
    @('1723857','Lettuce - Green Leaf','2','29.90','59.80'),
    @('1094721','Onion - Yellow Fresh','1','24.45','24.45'),
    @('0683696','Tomato - Fresh Sliced','14','31.71','443.94')
)

$r = 16
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
